$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '22.396.31'
$ws.Range('E2').Value = '  -0.11%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.573.60'
$ws.Range('E3').Value = '  +0.06%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '1.002'
$ws.Range('E5').Value = '  -0.04%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '291.19'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3765'
$ws.Range('E7').Value = '  +2.54%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '50.00'
$ws.Range('E8').Value = '  +1.23%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3416'
$ws.Range('E9').Value = '  +1.72%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.165'
$ws.Range('E10').Value = '  +0.04%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07677'
$ws.Range('E11').Value = '  +1.60%  '
$ws.Range('E12').Value = '  +0.09%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.36'
$ws.Range('E13').Value = '  +1.71%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.977'
$ws.Range('E14').Value = '  -1.59%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.914'
$ws.Range('E15').Value = '  +0.68%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.571.96'
$ws.Range('E16').Value = '  +0.13%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001138'
$ws.Range('E17').Value = '  +0.50%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '90.53'
$ws.Range('E18').Value = '  +1.35%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06740'
$ws.Range('E19').Value = '  -0.09%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.003'
$ws.Range('E20').Value = '  +0.04%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.78'
$ws.Range('E21').Value = '  +3.05%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.227'
$ws.Range('E22').Value = '  -0.25%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.5267'
$ws.Range('E23').Value = '  -3.84%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.00'
$ws.Range('E24').Value = '  +1.06%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '22.401.70'
$ws.Range('E25').Value = '  -0.09%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.405'
$ws.Range('E26').Value = '  +0.27%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.775'
$ws.Range('E27').Value = '  -5.82%  '
$ws.Range('E28').Value = '  +2.92%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '145.03'
$ws.Range('E29').Value = '  -0.37%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.068'
$ws.Range('E30').Value = '  +2.48%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '126.17'
$ws.Range('E31').Value = '  +1.21%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.746.77'
$ws.Range('E32').Value = '  +0.12%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.222'
$ws.Range('E33').Value = '  +0.00%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.017'
$ws.Range('E34').Value = '  +4.58%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.030'
$ws.Range('E35').Value = '  +1.64%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '10.05'
$ws.Range('E36').Value = '  -3.40%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.08564'
$ws.Range('E37').Value = '  +1.01%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02558'
$ws.Range('E38').Value = '  +1.98%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.2322'
$ws.Range('E39').Value = '  +1.52%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.06536'
$ws.Range('E40').Value = '  -0.32%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.331'
$ws.Range('E41').Value = '  +5.55%  '
$ws.Range('B42').Value = 'InternetComputer(DFINITY)'
$ws.Range('C42').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.483'
$ws.Range('E42').Value = '  +0.69%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '11.62'
$ws.Range('E43').Value = '  -1.38%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6466'
$ws.Range('E44').Value = '  +1.96%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '14.14'
$ws.Range('E45').Value = '  -2.74%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.002'
$ws.Range('E46').Value = '  +0.01%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.6032'
$ws.Range('E47').Value = '  +1.12%  '
$ws.Range('E48').Value = '  +0.25%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.300'
$ws.Range('E49').Value = '  +9.74%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.100'
$ws.Range('E50').Value = '  -0.44%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '125.53'
$ws.Range('E51').Value = '  +3.40%  '
